$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 - ASPM, trend period 5
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = 0.768783636774762
$ws.Range("K2").Value = 0.0060420140012664
$ws.Range("M2").Value = 0.0622805490331804
$ws.Range("N2").Value = 1.68771340817499
$ws.Range("P2").Value = "Likely improving"

# ---------------------------------------------------------------------------
# Row 3 - MCI, trend period 5
# ---------------------------------------------------------------------------
$ws.Range("J3").Value = 104.76
$ws.Range("K3").Value = -1.45219991245158
$ws.Range("M3").Value = 7.46439356934571
$ws.Range("N3").Value = -1.38621602944977

# ---------------------------------------------------------------------------
# Row 4 - QMCI, trend period 5
# ---------------------------------------------------------------------------
$ws.Range("F4").Value = 0.40324797025367
$ws.Range("J4").Value = 3.69
$ws.Range("K4").Value = -0.116188715612467
$ws.Range("M4").Value = 0.300872933065808
$ws.Range("N4").Value = -3.14874568055467

# ---------------------------------------------------------------------------
# Row 5 - ASPM, trend period 10
# ---------------------------------------------------------------------------
$ws.Range("F5").Value = 0.07620314197837499
$ws.Range("J5").Value = 0.3905
$ws.Range("K5").Value = -0.0177257347915243
$ws.Range("L5").Value = -0.0600818501232573
$ws.Range("M5").Value = -0.0000593300868559543
$ws.Range("N5").Value = -4.53924066364258
$ws.Range("P5").Value = "Very unlikely improving"

# ---------------------------------------------------------------------------
# Row 6 - MCI, trend period 10
# ---------------------------------------------------------------------------
$ws.Range("F6").Value = 0.185546684761349
$ws.Range("J6").Value = 105.38
$ws.Range("K6").Value = -1.24495446660885
$ws.Range("L6").Value = -3.93742015688464
$ws.Range("M6").Value = 1.11629870156652
$ws.Range("N6").Value = -1.18139539439063
$ws.Range("P6").Value = "Unlikely improving"

# ---------------------------------------------------------------------------
# Row 7 - QMCI, trend period 10
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 0.105248850027669
$ws.Range("J7").Value = 3.8845
$ws.Range("K7").Value = -0.174418245804006
$ws.Range("L7").Value = -0.488176992479633
$ws.Range("M7").Value = 0.0357109221588123
$ws.Range("N7").Value = -4.49010801400454
$ws.Range("P7").Value = "Unlikely improving"

# ---------------------------------------------------------------------------
# Row 8 - trend period 15 (parameter changes from MCI to ASPM)
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "ASPM (Macroinvertebrate Average Score Per Metric)"
$ws.Range("F8").Value = 0.008003077576582
$ws.Range("H8").Value = 1
$ws.Range("J8").Value = 0.427
$ws.Range("K8").Value = -0.0260751757706868
$ws.Range("L8").Value = -0.0426278832272651
$ws.Range("M8").Value = -0.0104270045003167
$ws.Range("N8").Value = -6.10659854114446
$ws.Range("P8").Value = "Exceptionally unlikely improving"

# ---------------------------------------------------------------------------
# New row 9 - MCI, trend period 15
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "Manganui o te Ao at Ashworth"
$ws.Range("B9").Value = "MCI (Macroinvertebrate Community Index)"
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = $false
$ws.Range("E9").Value = "ok"
$ws.Range("F9").Value = 0.012885200022937
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0.933333333333333
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 109
$ws.Range("K9").Value = -1.77118912080961
$ws.Range("L9").Value = -2.66500611382445
$ws.Range("M9").Value = -0.595554306388536
$ws.Range("N9").Value = -1.62494414753176
$ws.Range("O9").Value = "RepSite"
$ws.Range("P9").Value = "Extremely unlikely improving"
$ws.Range("Q9").Value = 1789685
$ws.Range("R9").Value = 5646155
$ws.Range("S9").Value = "Ruapehu District"
$ws.Range("T9").Value = "Whanganui"
$ws.Range("U9").Value = "Pipiriki"
$ws.Range("V9").Value = "Whai_5i"
$ws.Range("W9").Value = ""

# ---------------------------------------------------------------------------
# New row 10 - QMCI, trend period 15
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Manganui o te Ao at Ashworth"
$ws.Range("B10").Value = "QMCI (Quantitative Macroinvertebrate Community Index)"
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = $false
$ws.Range("E10").Value = "ok"
$ws.Range("F10").Value = 0.010744369559839
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 4.235
$ws.Range("K10").Value = -0.208204397394137
$ws.Range("L10").Value = -0.39308029610764
$ws.Range("M10").Value = -0.0916150120866006
$ws.Range("N10").Value = -4.91627856892885
$ws.Range("O10").Value = "RepSite"
$ws.Range("P10").Value = "Extremely unlikely improving"
$ws.Range("Q10").Value = 1789685
$ws.Range("R10").Value = 5646155
$ws.Range("S10").Value = "Ruapehu District"
$ws.Range("T10").Value = "Whanganui"
$ws.Range("U10").Value = "Pipiriki"
$ws.Range("V10").Value = "Whai_5i"
$ws.Range("W10").Value = ""
